$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps")
$ws.Range("G3").Value = "(Vx=5.0, Vy=3.5, Vz=0.8)"
Write-Host $ws.Range("G3").Value
